$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.920.67"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").Value = "1.895.18"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7768"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "243.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3121"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07355"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.51%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08080"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7704"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.495"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.83%  "
$ws.Range("D14").Value = "1.918.71"
$ws.Range("E14").Value = "  +1.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.223"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.98%  "
$ws.Range("D17").Value = "29.896.78"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.95"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007818"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9997"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.140.02"
$ws.Range("E22").Value = "  +0.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.107"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1582"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.431"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "164.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("E29").Value = "  -1.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.433"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.87%  "
$ws.Range("E31").Value = "  +0.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.479"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05565"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.28%  "
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.239"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.85%  "
$ws.Range("E36").Value = "  +1.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.004"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01931"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.798"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4466"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.11"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.27%  "
$ws.Range("D43").Value = "1.103.43"
$ws.Range("E43").Value = "  +7.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.962"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8507"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("E48").Value = "  +0.09%  "
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.511"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.737"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.031"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.80%  "
